$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "40"
    6  = "0.00064"
    7  = "0.00027"
    8  = "0.00006"
    9  = "0.00049"
    10 = "0.00051"
    11 = "0.00059"
    12 = "0.01081"
    44 = "99.94"
    45 = "0.01"
    46 = "19"
}

foreach ($rowIndex in $updates.Keys) {
    $t.Rows.Item($rowIndex).Cells.Item(1).Range.Text = $updates[$rowIndex]
}
